# Update "Generate Report for Handback" timestamps.
#
# Mapping of cells -> old / new timestamp values (from the shared-string diff):
#   Overview!G2  ("Latest HO Xliff Generate Date")      2016-08-26 07:03:21 -> 2016-08-26 07:04:08
#   zh-cn!H2     ("Correspond Handoff Datetime")         2016-08-26 07:03:16 -> 2016-08-26 07:03:58
#   zh-cn!K2     ("Correspond Handback DateTime")        2016-08-26 07:03:32 -> 2016-08-26 07:04:32
#   de-de!H2     ("Correspond Handoff Datetime")         2016-08-26 07:03:21 -> 2016-08-26 07:04:08 (shared string w/ Overview!G2)
#   de-de!K2     ("Correspond Handback DateTime")        2016-08-26 07:03:38 -> 2016-08-26 07:04:39

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-26 07:04:08"

$wsZhCn.Range("H2").Value = "2016-08-26 07:03:58"
$wsZhCn.Range("K2").Value = "2016-08-26 07:04:32"

$wsDeDe.Range("H2").Value = "2016-08-26 07:04:08"
$wsDeDe.Range("K2").Value = "2016-08-26 07:04:39"
